{"js": "// Adds the \"Extended Project Information\" section (Impact Statement,\n// Sustainability Plan, Appendices A-F) after the existing\n// \"Research Network Leadership:\" paragraph, matching the target diff.\n\nconst PARAGRAPHS = [\n  [\n    [\n      \"Extended Project Information\",\n      true\n    ]\n  ],\n  [\n    [\n      \"Project Statistics:\",\n      true\n    ],\n    [\n      \" Six core publications with 28+ total citations received, 12 Zenodo deposits ensuring reproducibility, 8 conference presentations across 8 countries, 2 PhD researchers trained, and CHF 90,000 in additional funding secured.\",\n      false\n    ]\n  ],\n  [\n    [\n      \"1.4 Impact Statement\",\n      true\n    ]\n  ],\n  [\n    [\n      \"Scientific Impact:\",\n      true\n    ],\n    [\n      \" Our publications have received 28+ citations within 12 months of publication, indicating rapid adoption by the research community. The two-step machine learning methodology combining network centrality with traditional credit features has been referenced in subsequent studies on P2P lending risk assessment, and the systematic literature review provides a foundational reference for researchers entering the field of graph-based credit modeling.\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Economic Impact:\",\n      true\n    ],\n    [\n      \" The developed models and open-source code are directly applicable to P2P lending platforms for improved credit risk assessment. By enabling more accurate default prediction, platforms can better price loans to reflect true risk, reduce losses from defaults, and offer more competitive rates to creditworthy borrowers. The interpretability framework addresses regulatory requirements, reducing compliance costs for platforms adopting automated credit decisions.\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Social Impact:\",\n      true\n    ],\n    [\n      \" Improved credit risk models contribute to financial inclusion by enabling P2P platforms to serve borrowers who may be underserved by traditional banking. More accurate risk assessment reduces adverse selection problems, protecting retail investors who fund P2P loans from excessive default losses. The transparency framework enhances borrower trust in automated credit decisions by providing explanations for lending outcomes.\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Policy Impact:\",\n      true\n    ],\n    [\n      \" The PI's leadership of COST Action CA19130 facilitated policy discussions at EU level, including events in Brussels addressing AI in finance policy implications and contributing to the broader discourse on responsible AI adoption in financial services.\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Educational Impact:\",\n      true\n    ],\n    [\n      \" The project trained two PhD researchers in cutting-edge methods at the intersection of network science, machine learning, and finance. The publication \\\"Towards a new PhD Curriculum for Digital Finance\\\" (Open Research Europe, 2024, DOI: 10.12688/openreseurope.16513.1) disseminates best practices for doctoral training in this emerging field, contributing to curriculum development beyond this specific project.\",\n      false\n    ]\n  ],\n  [\n    [\n      \"1.5 Sustainability Plan\",\n      true\n    ]\n  ],\n  [\n    [\n      \"Data Preservation:\",\n      true\n    ],\n    [\n      \" Twelve Zenodo deposits are archived with persistent DOIs ensuring permanent accessibility and citability through CERN's infrastructure. The curated Bondora P2P lending dataset is archived at the Open Science Framework (OSF). Code repositories are maintained under the Digital-AI-Finance organization on GitHub. All outputs are released under Creative Commons Attribution 4.0 (CC-BY 4.0) licensing, enabling unrestricted reuse with attribution.\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Code Maintainability:\",\n      true\n    ],\n    [\n      \" All Jupyter notebooks and Python/R scripts include dependency specifications (requirements.txt, environment files) enabling reproduction with specified package versions. Reproducibility has been verified through independent testing. Documentation is embedded in code through comments and supplementary README files.\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Knowledge Transfer Continuation:\",\n      true\n    ],\n    [\n      \" COST Action CA19130 (Fintech and AI in Finance) continues beyond project end with the PI serving as Action Chair. The MSCA Industrial Doctoral Network on Digital Finance continues training next-generation researchers with the PI as Coordinator. Digital finance research continues at Bern University of Applied Sciences building on this project's foundations.\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Appendix A:\",\n      true\n    ],\n    [\n      \" Peer-Reviewed Publications\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Liu, Y., Baals, L.J., Osterrieder, J., Hadji-Misheva, B. (2024). Leveraging network topology for credit risk assessment in P2P lending: A comparative study under the lens of machine learning. Expert Systems with Applications, 252(B), 124100. DOI: 10.1016/j.eswa.2024.124100. 17 citations.\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Liu, Y., Baals, L.J., Osterrieder, J., Hadji-Misheva, B. (2024). Network centrality and credit risk: A comprehensive analysis of peer-to-peer lending dynamics. Finance Research Letters, 63, 105308. DOI: 10.1016/j.frl.2024.105308. 11 citations.\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Baumohl, E., Lyocsa, S., Vasanicova, P. (2024). Macroeconomic environment and the future performance of loans: Evidence from three peer-to-peer platforms. International Review of Financial Analysis, 95, 103416. DOI: 10.1016/j.irfa.2024.103416.\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Baals, L.J., Osterrieder, J., Hadji-Misheva, B., Liu, Y. (2024). Towards a new PhD Curriculum for Digital Finance. Open Research Europe, 4, 16513. DOI: 10.12688/openreseurope.16513.1.\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Submitted: Baals, L.J., et al. (2025). Network Evidence on Credit-Risk Pricing in P2P Lending. SSRN 5276337. Baals, L.J., et al. (2025). State-Dependent Pricing in FinTech Credit: Evidence from P2P Lending. SSRN 5421207.\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Appendix B:\",\n      true\n    ],\n    [\n      \" Open Science Deposits (Zenodo)\",\n      false\n    ]\n  ],\n  [\n    [\n      \"1. COST FinAI Meets Istanbul Conference Event May 20-21, 2024. Baals, Lennart John (2024). Presentation. https://zenodo.org/records/17964900\",\n      false\n    ]\n  ],\n  [\n    [\n      \"2. State-Dependent Pricing in FinTech Credit: Evidence from P2P Lending. Baals, Lennart John (2025). Working paper. https://zenodo.org/records/17990398\",\n      false\n    ]\n  ],\n  [\n    [\n      \"3. Network Evidence on Credit-Risk Pricing in P2P Lending. Baals, Lennart John (2025). Working paper. https://zenodo.org/records/17990873\",\n      false\n    ]\n  ],\n  [\n    [\n      \"4. Leveraging Network Topology for Credit Risk Assessment in P2P Lending: A Comparative Study under the Lens of Machine Learning. Baals, Lennart John (2025). Journal article. https://zenodo.org/records/17991107\",\n      false\n    ]\n  ],\n  [\n    [\n      \"5. PhD Qualifier Report and Presentation delivered by Lennart John Baals at the University of Twente. Baals, Lennart John (2025). Proposal. https://zenodo.org/records/17992215\",\n      false\n    ]\n  ],\n  [\n    [\n      \"6. A Systematic Literature Review on Graph-Based Models in Credit Risk Assessment. Baals, Lennart John (2025). Presentation. https://zenodo.org/records/17992322\",\n      false\n    ]\n  ],\n  [\n    [\n      \"7. Leveraging Network Topology for Credit Risk Assessment in P2P Lending (Bern Conference 2023). Baals, Lennart John (2025). Presentation. https://zenodo.org/records/17992484\",\n      false\n    ]\n  ],\n  [\n    [\n      \"8. Identifying Mispriced Loans through Interest Rate-Based Network Analysis and Clustering in P2P Lending Markets. Baals, Lennart John (2025). Presentation. https://zenodo.org/records/17992591\",\n      false\n    ]\n  ],\n  [\n    [\n      \"9. Data and Code to reproduce results in paper \\\"Network centrality and credit risk\\\". Liu, Yiting (2024). Computational notebook. https://zenodo.org/records/17989119\",\n      false\n    ]\n  ],\n  [\n    [\n      \"10. Data and Code to reproduce results in paper \\\"Credit Risk Prediction via Graph Neural Networks with Homophily-Guided Graph Construction\\\". Liu, Yiting (2026). Computational notebook. https://zenodo.org/records/17990002\",\n      false\n    ]\n  ],\n  [\n    [\n      \"11. Data and Code to reproduce results in paper \\\"Explaining Regime Dynamics: A Tree-based Interpretation Framework for R2-RD Models\\\". Liu, Yiting (2026). Computational notebook. https://zenodo.org/records/17990140\",\n      false\n    ]\n  ],\n  [\n    [\n      \"12. Data and Code to reproduce results in paper \\\"Leveraging network topology for credit risk assessment in P2P lending\\\". Liu, Yiting (2024). Computational notebook. https://zenodo.org/records/17990581\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Appendix C:\",\n      true\n    ],\n    [\n      \" Academic Events\",\n      false\n    ]\n  ],\n  [\n    [\n      \"December 2024: 4th International Symposium on Big Data and AI, Hong Kong (Systematic Literature Review on Graph-Based Credit Models). September 2024: 8th Bern Conference on Fintech and AI in Finance, Switzerland. September 2024: AI Finance Insights: Pioneering the Future of Fintech, Istanbul. May 2024: COST FinAI Meets Istanbul Conference, Turkey. December 2023: 16th ERCIM Conference on Computational and Methodological Statistics, Berlin. September 2023: 8th European COST Conference on AI in Finance, Bern. September 2023: European Summer School in Financial Mathematics, Delft.\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Appendix D:\",\n      true\n    ],\n    [\n      \" Dataset\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Bondora P2P Lending Dataset. Coverage: June 2009 - April 2022. Sample: 231,039 borrowers, 112 variables. DOI: 10.21227/33kz-0s65. License: CC-BY 4.0.\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Appendix E:\",\n      true\n    ],\n    [\n      \" International Collaborations\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Masaryk University (Czech Republic), Columbia University (USA), American University of Sharjah (UAE), Renmin University of China (China), University of Manchester (UK).\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Appendix F:\",\n      true\n    ],\n    [\n      \" PhD Researchers\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Lennart John Baals: PhD In Progress, BFH/University of Twente, Graph-based credit models and network analysis for credit risk assessment.\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Yiting Liu: PhD In Progress, BFH/University of Twente, P2P lending risk modeling and network topology for credit risk.\",\n      false\n    ]\n  ],\n  [\n    [\n      \"Report submitted to:\",\n      true\n    ],\n    [\n      \" Swiss National Science Foundation (SNSF). Report date: December 2025. Data source: https://data.snf.ch/grants/grant/205487\",\n      false\n    ]\n  ]\n];\n\n// Locate the anchor paragraph: the last non-empty paragraph in the body,\n// which ends with the \"Research Network Leadership:\" text (immediately\n// before the final empty paragraph / section break).\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nlet anchor = null;\nfor (let i = paras.items.length - 1; i >= 0; i--) {\n  if (paras.items[i].text && paras.items[i].text.indexOf(\"Research Network Leadership:\") !== -1) {\n    anchor = paras.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find anchor paragraph 'Research Network Leadership:'\");\n}\n\n// Insert each new paragraph, in order, directly after the anchor so the\n// final order matches the source document.\nlet current = anchor;\nfor (const runs of PARAGRAPHS) {\n  const firstText = runs[0][0];\n  const firstBold = runs[0][1];\n  const newPara = current.insertParagraph(firstText, Word.InsertLocation.after);\n  newPara.font.bold = firstBold;\n  for (let i = 1; i < runs.length; i++) {\n    const [text, bold] = runs[i];\n    const range = newPara.insertText(text, Word.InsertLocation.end);\n    range.font.bold = bold;\n  }\n  current = newPara;\n}\n\nawait context.sync();\n", "ps1": "# Adds the \"Extended Project Information\" section (Project Statistics,\n# 1.4 Impact Statement, 1.5 Sustainability Plan, Appendices A-F) after the\n# existing \"Research Network Leadership:\" paragraph, matching the target diff.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph: the one whose text contains\n# \"Research Network Leadership:\" -- the last substantive paragraph of the\n# existing section, immediately before the trailing empty paragraph.\n$anchorIdx = -1\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $cand = $d.Paragraphs.Item($i)\n    if ($cand.Range.Text -like '*Research Network Leadership:*') {\n        $anchorIdx = $i\n        break\n    }\n}\nif ($anchorIdx -eq -1) {\n    throw \"Could not find anchor paragraph 'Research Network Leadership:'\"\n}\n\n# Build the list of new paragraphs. Each paragraph is an array of\n# (text, isBold) run pairs, inserted in order directly after the anchor.\n# The leading comma on each $pN keeps PowerShell from unwrapping\n# single-run paragraphs when the pieces are concatenated below.\n\n$p0 = @(,@('Extended Project Information', $true))\n$p1 = @(@('Project Statistics:', $true), @(' Six core publications with 28+ total citations received, 12 Zenodo deposits ensuring reproducibility, 8 conference presentations across 8 countries, 2 PhD researchers trained, and CHF 90,000 in additional funding secured.', $false))\n$p2 = @(,@('1.4 Impact Statement', $true))\n$p3 = @(@('Scientific Impact:', $true), @(' Our publications have received 28+ citations within 12 months of publication, indicating rapid adoption by the research community. The two-step machine learning methodology combining network centrality with traditional credit features has been referenced in subsequent studies on P2P lending risk assessment, and the systematic literature review provides a foundational reference for researchers entering the field of graph-based credit modeling.', $false))\n$p4 = @(@('Economic Impact:', $true), @(' The developed models and open-source code are directly applicable to P2P lending platforms for improved credit risk assessment. By enabling more accurate default prediction, platforms can better price loans to reflect true risk, reduce losses from defaults, and offer more competitive rates to creditworthy borrowers. The interpretability framework addresses regulatory requirements, reducing compliance costs for platforms adopting automated credit decisions.', $false))\n$p5 = @(@('Social Impact:', $true), @(' Improved credit risk models contribute to financial inclusion by enabling P2P platforms to serve borrowers who may be underserved by traditional banking. More accurate risk assessment reduces adverse selection problems, protecting retail investors who fund P2P loans from excessive default losses. The transparency framework enhances borrower trust in automated credit decisions by providing explanations for lending outcomes.', $false))\n$p6 = @(@('Policy Impact:', $true), @(' The PI''s leadership of COST Action CA19130 facilitated policy discussions at EU level, including events in Brussels addressing AI in finance policy implications and contributing to the broader discourse on responsible AI adoption in financial services.', $false))\n$p7 = @(@('Educational Impact:', $true), @(' The project trained two PhD researchers in cutting-edge methods at the intersection of network science, machine learning, and finance. The publication \"Towards a new PhD Curriculum for Digital Finance\" (Open Research Europe, 2024, DOI: 10.12688/openreseurope.16513.1) disseminates best practices for doctoral training in this emerging field, contributing to curriculum development beyond this specific project.', $false))\n$p8 = @(,@('1.5 Sustainability Plan', $true))\n$p9 = @(@('Data Preservation:', $true), @(' Twelve Zenodo deposits are archived with persistent DOIs ensuring permanent accessibility and citability through CERN''s infrastructure. The curated Bondora P2P lending dataset is archived at the Open Science Framework (OSF). Code repositories are maintained under the Digital-AI-Finance organization on GitHub. All outputs are released under Creative Commons Attribution 4.0 (CC-BY 4.0) licensing, enabling unrestricted reuse with attribution.', $false))\n$p10 = @(@('Code Maintainability:', $true), @(' All Jupyter notebooks and Python/R scripts include dependency specifications (requirements.txt, environment files) enabling reproduction with specified package versions. Reproducibility has been verified through independent testing. Documentation is embedded in code through comments and supplementary README files.', $false))\n$p11 = @(@('Knowledge Transfer Continuation:', $true), @(' COST Action CA19130 (Fintech and AI in Finance) continues beyond project end with the PI serving as Action Chair. The MSCA Industrial Doctoral Network on Digital Finance continues training next-generation researchers with the PI as Coordinator. Digital finance research continues at Bern University of Applied Sciences building on this project''s foundations.', $false))\n$p12 = @(@('Appendix A:', $true), @(' Peer-Reviewed Publications', $false))\n$p13 = @(,@('Liu, Y., Baals, L.J., Osterrieder, J., Hadji-Misheva, B. (2024). Leveraging network topology for credit risk assessment in P2P lending: A comparative study under the lens of machine learning. Expert Systems with Applications, 252(B), 124100. DOI: 10.1016/j.eswa.2024.124100. 17 citations.', $false))\n$p14 = @(,@('Liu, Y., Baals, L.J., Osterrieder, J., Hadji-Misheva, B. (2024). Network centrality and credit risk: A comprehensive analysis of peer-to-peer lending dynamics. Finance Research Letters, 63, 105308. DOI: 10.1016/j.frl.2024.105308. 11 citations.', $false))\n$p15 = @(,@('Baumohl, E., Lyocsa, S., Vasanicova, P. (2024). Macroeconomic environment and the future performance of loans: Evidence from three peer-to-peer platforms. International Review of Financial Analysis, 95, 103416. DOI: 10.1016/j.irfa.2024.103416.', $false))\n$p16 = @(,@('Baals, L.J., Osterrieder, J., Hadji-Misheva, B., Liu, Y. (2024). Towards a new PhD Curriculum for Digital Finance. Open Research Europe, 4, 16513. DOI: 10.12688/openreseurope.16513.1.', $false))\n$p17 = @(,@('Submitted: Baals, L.J., et al. (2025). Network Evidence on Credit-Risk Pricing in P2P Lending. SSRN 5276337. Baals, L.J., et al. (2025). State-Dependent Pricing in FinTech Credit: Evidence from P2P Lending. SSRN 5421207.', $false))\n$p18 = @(@('Appendix B:', $true), @(' Open Science Deposits (Zenodo)', $false))\n$p19 = @(,@('1. COST FinAI Meets Istanbul Conference Event May 20-21, 2024. Baals, Lennart John (2024). Presentation. https://zenodo.org/records/17964900', $false))\n$p20 = @(,@('2. State-Dependent Pricing in FinTech Credit: Evidence from P2P Lending. Baals, Lennart John (2025). Working paper. https://zenodo.org/records/17990398', $false))\n$p21 = @(,@('3. Network Evidence on Credit-Risk Pricing in P2P Lending. Baals, Lennart John (2025). Working paper. https://zenodo.org/records/17990873', $false))\n$p22 = @(,@('4. Leveraging Network Topology for Credit Risk Assessment in P2P Lending: A Comparative Study under the Lens of Machine Learning. Baals, Lennart John (2025). Journal article. https://zenodo.org/records/17991107', $false))\n$p23 = @(,@('5. PhD Qualifier Report and Presentation delivered by Lennart John Baals at the University of Twente. Baals, Lennart John (2025). Proposal. https://zenodo.org/records/17992215', $false))\n$p24 = @(,@('6. A Systematic Literature Review on Graph-Based Models in Credit Risk Assessment. Baals, Lennart John (2025). Presentation. https://zenodo.org/records/17992322', $false))\n$p25 = @(,@('7. Leveraging Network Topology for Credit Risk Assessment in P2P Lending (Bern Conference 2023). Baals, Lennart John (2025). Presentation. https://zenodo.org/records/17992484', $false))\n$p26 = @(,@('8. Identifying Mispriced Loans through Interest Rate-Based Network Analysis and Clustering in P2P Lending Markets. Baals, Lennart John (2025). Presentation. https://zenodo.org/records/17992591', $false))\n$p27 = @(,@('9. Data and Code to reproduce results in paper \"Network centrality and credit risk\". Liu, Yiting (2024). Computational notebook. https://zenodo.org/records/17989119', $false))\n$p28 = @(,@('10. Data and Code to reproduce results in paper \"Credit Risk Prediction via Graph Neural Networks with Homophily-Guided Graph Construction\". Liu, Yiting (2026). Computational notebook. https://zenodo.org/records/17990002', $false))\n$p29 = @(,@('11. Data and Code to reproduce results in paper \"Explaining Regime Dynamics: A Tree-based Interpretation Framework for R2-RD Models\". Liu, Yiting (2026). Computational notebook. https://zenodo.org/records/17990140', $false))\n$p30 = @(,@('12. Data and Code to reproduce results in paper \"Leveraging network topology for credit risk assessment in P2P lending\". Liu, Yiting (2024). Computational notebook. https://zenodo.org/records/17990581', $false))\n$p31 = @(@('Appendix C:', $true), @(' Academic Events', $false))\n$p32 = @(,@('December 2024: 4th International Symposium on Big Data and AI, Hong Kong (Systematic Literature Review on Graph-Based Credit Models). September 2024: 8th Bern Conference on Fintech and AI in Finance, Switzerland. September 2024: AI Finance Insights: Pioneering the Future of Fintech, Istanbul. May 2024: COST FinAI Meets Istanbul Conference, Turkey. December 2023: 16th ERCIM Conference on Computational and Methodological Statistics, Berlin. September 2023: 8th European COST Conference on AI in Finance, Bern. September 2023: European Summer School in Financial Mathematics, Delft.', $false))\n$p33 = @(@('Appendix D:', $true), @(' Dataset', $false))\n$p34 = @(,@('Bondora P2P Lending Dataset. Coverage: June 2009 - April 2022. Sample: 231,039 borrowers, 112 variables. DOI: 10.21227/33kz-0s65. License: CC-BY 4.0.', $false))\n$p35 = @(@('Appendix E:', $true), @(' International Collaborations', $false))\n$p36 = @(,@('Masaryk University (Czech Republic), Columbia University (USA), American University of Sharjah (UAE), Renmin University of China (China), University of Manchester (UK).', $false))\n$p37 = @(@('Appendix F:', $true), @(' PhD Researchers', $false))\n$p38 = @(,@('Lennart John Baals: PhD In Progress, BFH/University of Twente, Graph-based credit models and network analysis for credit risk assessment.', $false))\n$p39 = @(,@('Yiting Liu: PhD In Progress, BFH/University of Twente, P2P lending risk modeling and network topology for credit risk.', $false))\n$p40 = @(@('Report submitted to:', $true), @(' Swiss National Science Foundation (SNSF). Report date: December 2025. Data source: https://data.snf.ch/grants/grant/205487', $false))\n\n$paragraphs = @(,$p0) + @(,$p1) + @(,$p2) + @(,$p3) + @(,$p4) + @(,$p5) + @(,$p6) + @(,$p7) + @(,$p8) + @(,$p9) + @(,$p10) + @(,$p11) + @(,$p12) + @(,$p13) + @(,$p14) + @(,$p15) + @(,$p16) + @(,$p17) + @(,$p18) + @(,$p19) + @(,$p20) + @(,$p21) + @(,$p22) + @(,$p23) + @(,$p24) + @(,$p25) + @(,$p26) + @(,$p27) + @(,$p28) + @(,$p29) + @(,$p30) + @(,$p31) + @(,$p32) + @(,$p33) + @(,$p34) + @(,$p35) + @(,$p36) + @(,$p37) + @(,$p38) + @(,$p39) + @(,$p40)\n\n# Walk forward from the anchor paragraph, inserting one new empty\n# paragraph at a time and filling it with its runs (text + bold flag).\n$currentIdx = $anchorIdx\nforeach ($runs in $paragraphs) {\n    $cur = $d.Paragraphs.Item($currentIdx)\n    $r = $cur.Range\n    $r.Collapse(0) | Out-Null   # wdCollapseEnd = 0\n    $r.InsertParagraphAfter()\n    $currentIdx = $currentIdx + 1\n    $newPara = $d.Paragraphs.Item($currentIdx)\n    $cursor = $newPara.Range.Start\n    foreach ($run in $runs) {\n        $text = $run[0]\n        $bold = $run[1]\n        $ins = $d.Range($cursor, $cursor)\n        $ins.InsertAfter($text)\n        $runRange = $d.Range($cursor, $cursor + $text.Length)\n        if ($bold) {\n            $runRange.Font.Bold = 1\n        } else {\n            $runRange.Font.Bold = 0\n        }\n        $cursor = $cursor + $text.Length\n    }\n}\n\n"}
